# Project Sample Project is saved (SAVE): cell B11 on sheet "Rules" changes
# from the rule-id text "R40" to the text "1" (rule row 40 -> 1), keeping
# its existing cell style/formatting untouched.
#
# A plain `.Value = "1"` would be auto-coerced to the *number* 1 (and any
# NumberFormat/quote-prefix trick to force text re-keys the cell style).
# Instead, compute the text via TEXT() and paste back as a value only, so
# the destination keeps its original style and becomes a genuine text
# (shared-string) cell, exactly like typing ="1" into the cell and
# converting the formula result to a static value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
$cell.Formula = "=TEXT(1,0)"
$cell.Copy()
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
$excel.CutCopyMode = $false
